$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Kunlik hujjat topshirgan" (E4:E13) column contents that were
# previously filled in, leaving the cell formatting untouched.
$ws.Range("E4:E13").ClearContents() | Out-Null

# Reflect the user's selection after clearing the range (active cell E4,
# selection spanning E4:E13).
$ws.Range("E4:E13").Select() | Out-Null
